$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 updates
$ws.Range("C6").Value = 2.01
$ws.Range("D6").Value = 25
$ws.Range("E6").Value = 3.72
$ws.Range("J6").Value = 5
$ws.Range("M6").Value = 5
$ws.Range("Q6").Value = 4
$ws.Range("S6").Value = 2.01
$ws.Range("T6").Value = 25
$ws.Range("U6").Value = 3.72

# Row 9 updates
$ws.Range("M9").Value = 16
$ws.Range("S9").Value = 1.71
$ws.Range("T9").Value = 80
$ws.Range("U9").Value = 3.35
